# Update "想去人数" (interested-count) figures on the 展览 (Exhibitions)
# and 全部类型 (All Types) sheets to match the freshly scraped stats.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 sheet — column F values
$wsExhibition.Range("F5").Value  = 381
$wsExhibition.Range("F6").Value  = 668
$wsExhibition.Range("F8").Value  = 2104
$wsExhibition.Range("F9").Value  = 6
$wsExhibition.Range("F10").Value = 10928
$wsExhibition.Range("F11").Value = 184
$wsExhibition.Range("F15").Value = 10726
$wsExhibition.Range("F16").Value = 431
$wsExhibition.Range("F18").Value = 6
$wsExhibition.Range("F19").Value = 748
$wsExhibition.Range("F20").Value = 5368
$wsExhibition.Range("F22").Value = 3388

# 全部类型 sheet — column F values (same events, different row offsets)
$wsAllTypes.Range("F5").Value  = 381
$wsAllTypes.Range("F6").Value  = 668
$wsAllTypes.Range("F9").Value  = 2104
$wsAllTypes.Range("F11").Value = 6
$wsAllTypes.Range("F13").Value = 10928
$wsAllTypes.Range("F14").Value = 184
$wsAllTypes.Range("F18").Value = 10726
$wsAllTypes.Range("F19").Value = 431
$wsAllTypes.Range("F21").Value = 6
$wsAllTypes.Range("F22").Value = 748
$wsAllTypes.Range("F23").Value = 5368
$wsAllTypes.Range("F25").Value = 3388
